$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 currently holds the "old" record (date 44399, etc.).
# We need to:
#   1) Copy the existing row 20 down to row 21 (so the old record is preserved there)
#   2) Update row 20 in place with the new record's values

$ws.Range("A20:R20").Copy($ws.Range("A21:R21")) | Out-Null
$excel.CutCopyMode = 0

# Now overwrite row 20 with the updated values
$ws.Range("D20").Value = 44448
$ws.Range("J20").Value = 45
$ws.Range("K20").Value = 32000
$ws.Range("L20").Value = 32000
$ws.Range("M20").Value = 32000
$ws.Range("P20").Value = 1280
